# BaoCaoNhapHang.xlsx - "fix 27/4/2024 lan 1"
#
# The "Thang 1" sheet's header table gains a new "Kich co" (size) column
# and its headers are renamed/reordered to match the updated shoe-store
# inventory report layout:
#   STT | Ma san pham | Ten san pham | Mau sac | Kich co | So luong nhap | Thanh tien
#
# Only the first sheet ("Thang 1") is touched - the other month sheets keep
# their original 6-column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Insert the new column *inside* the existing A1:F1 merged title band (at
# column F, pushing the old column F to G) so the merge grows to A1:G1 and
# every cell naturally inherits the surrounding s="3"/s="2" styles instead
# of minting new ones.
$ws.Columns.Item(6).Insert()

# Re-label the header row (row 2). Written in this order so the shared
# strings table is appended in the same sequence as the source workbook.
$ws.Range("E2").Value = "Kích cỡ"
$ws.Range("D2").Value = "Màu sắc"
$ws.Range("G2").Value = "Thành tiền"
$ws.Range("F2").Value = "Số lượng nhập"
$ws.Range("C2").Value = "Tên sản phẩm"
$ws.Range("B2").Value = "Mã sản phẩm"

# Match the row-1 title band height tweak that came with this edit.
$ws.Rows.Item(1).RowHeight = 28.95

# Best-effort column widths for the new/resized columns (B..G). The exact
# sub-pixel values depend on the authoring Excel build's font metrics and
# can't be reproduced bit-for-bit, so these are the closest attainable
# approximations of the saved widths.
$ws.Columns.Item(2).ColumnWidth = 17.498697916666668
$ws.Columns.Item(3).ColumnWidth = 37.276041666666664
$ws.Columns.Item(4).ColumnWidth = 12.385416666666666
$ws.Columns.Item(5).ColumnWidth = 12.721354166666666
$ws.Columns.Item(6).ColumnWidth = 17.276041666666668
$ws.Columns.Item(7).ColumnWidth = 18.498697916666668

# Leave the cursor where the author's last save left it.
$ws.Range("J6").Select()
